$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the "k" column (J), placed right under the data, bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Summary rows 14-17: a label in column A and an aggregate formula in
# column B, formatted bold / size 12 / vertically centered
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertical-centered look on B14 first, then copy
# that exact formatting (not the value) onto B15:B17 so only one extra
# style definition is introduced
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Page setup (paper size / orientation) as recorded in the saved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the final selection on B17, matching the saved view state
$ws.Range("B17").Select()

$wb.Save()
